# Updated schedules, audio and instructions
# The two trial rows (row 2 / row 3) were re-ordered: the data that used
# to live on row 3 now lives on row 2, and vice-versa. Columns A (trial
# index), G (y_nrSteps), I (trialLength) and J (version) already matched
# between the two rows, so only B, C, D, E, F and H actually change value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 2 (previously row 3's values)
$ws.Range("B2").Value = 3
$ws.Range("C2").Value = 2
$ws.Range("D2").Value = 8
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 5
$ws.Range("H2").Value = 66

# New row 3 (previously row 2's values)
$ws.Range("B3").Value = 4
$ws.Range("C3").Value = 0
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("H3").Value = 16

# Match the saved selection state recorded in the workbook (cell J10).
$ws.Range("J10").Select()
